# Auto-generated Excel COM-interop script applying the diff changes
$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions): numeric "want-to-go" count updates ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 621
$wsExpo.Range("F3").Value = 5847
$wsExpo.Range("F5").Value = 463
$wsExpo.Range("F7").Value = 1017
$wsExpo.Range("F8").Value = 390
$wsExpo.Range("F11").Value = 3122
$wsExpo.Range("F12").Value = 1955
$wsExpo.Range("F15").Value = 200
$wsExpo.Range("F16").Value = 78
$wsExpo.Range("F17").Value = 171
$wsExpo.Range("F19").Value = 988
$wsExpo.Range("F22").Value = 57
$wsExpo.Range("F23").Value = 3660
$wsExpo.Range("F24").Value = 1165
$wsExpo.Range("F25").Value = 2903
$wsExpo.Range("F26").Value = 289
$wsExpo.Range("F27").Value = 2259
$wsExpo.Range("F28").Value = 4207
$wsExpo.Range("F30").Value = 930
$wsExpo.Range("F32").Value = 1340
$wsExpo.Range("F33").Value = 96
$wsExpo.Range("F34").Value = 3
$wsExpo.Range("F35").Value = 20
$wsExpo.Range("F36").Value = 23
$wsExpo.Range("F41").Value = 1099
$wsExpo.Range("F42").Value = 708
$wsExpo.Range("F43").Value = 597
$wsExpo.Range("F44").Value = 430
$wsExpo.Range("F46").Value = 97
$wsExpo.Range("F47").Value = 2
$wsExpo.Range("F49").Value = 3605

# ---- Sheet "演出" (performances): numeric "want-to-go" count updates ----
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 9
$wsShow.Range("F10").Value = 914
$wsShow.Range("F26").Value = 9

# ---- Sheet "全部类型" (all types): numeric "want-to-go" count updates ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 621
$wsAll.Range("F3").Value = 5847
$wsAll.Range("F6").Value = 9
$wsAll.Range("F7").Value = 390
$wsAll.Range("F9").Value = 3122
$wsAll.Range("F11").Value = 1955
$wsAll.Range("F14").Value = 200
$wsAll.Range("F15").Value = 914
$wsAll.Range("F29").Value = 930
$wsAll.Range("F30").Value = 1340
$wsAll.Range("F36").Value = 1099
$wsAll.Range("F38").Value = 708
$wsAll.Range("F40").Value = 430
$wsAll.Range("F45").Value = 97
$wsAll.Range("F46").Value = 9
$wsAll.Range("F48").Value = 3606

# ---- Sheet "全部类型": row content shift rows 17-27 (new exhibition inserted at row 17) ----
$wsAll.Range("B17").Value = "2024-06-23"
$wsAll.Range("C17").Value = "杭州·第二届白日梦次元动漫嘉年华"
$wsAll.Range("D17").Value = "康候圣街99号 顺丰创新中心"
$wsAll.Range("E17").Value = "2024.06.23 10:00-06.23 17:00"
$wsAll.Range("F17").Value = 78
$wsAll.Range("G17").Value = 68
$wsAll.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=86307"
$wsAll.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202405/qHcyIUL31715752173541.jpeg"

$wsAll.Range("B18").Value = "2024-06-23"
$wsAll.Range("C18").Value = "杭州·第五人格ONLY2.0"
$wsAll.Range("D18").Value = "十四号大街431号 江滨篮球馆"
$wsAll.Range("E18").Value = "2024.06.23 10:00-06.23 17:00"
$wsAll.Range("F18").Value = 171
$wsAll.Range("G18").Value = 60
$wsAll.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=85710"
$wsAll.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202405/ULUN091G1715762966375.jpeg"

$wsAll.Range("B19").Value = "2024-06-29"
$wsAll.Range("C19").Value = "杭州·乌托邦次元聚会3.0·二次元全女性夜场"
$wsAll.Range("D19").Value = "保淑路2号 The Queen皇后"
$wsAll.Range("E19").Value = "2024.06.29 13:00-06.29 19:00"
$wsAll.Range("F19").Value = 988
$wsAll.Range("G19").Value = 188
$wsAll.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=84558"
$wsAll.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202404/XyOkWYv31713435061841.jpeg"

$wsAll.Range("B20").Value = "2024-06-30"
$wsAll.Range("C20").Value = "杭州·热血番ONLY1.0"
$wsAll.Range("D20").Value = "康候圣街99号 顺丰创新中心"
$wsAll.Range("E20").Value = "2024.06.30 10:00-06.30 17:00"
$wsAll.Range("F20").Value = 360
$wsAll.Range("G20").Value = 68
$wsAll.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=85042"
$wsAll.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202404/HSZzsbLs1714221391790.jpeg"

$wsAll.Range("B21").Value = "2024-07-13"
$wsAll.Range("C21").Value = "杭州·AD04动漫展"
$wsAll.Range("D21").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$wsAll.Range("E21").Value = "2024.07.13 10:00-07.14 17:00"
$wsAll.Range("F21").Value = 3660
$wsAll.Range("G21").Value = 75
$wsAll.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=85012"
$wsAll.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202405/y1iKqqnh1715326769523.jpeg"

$wsAll.Range("B22").Value = "2024-07-13"
$wsAll.Range("C22").Value = "杭州·【早鸟6折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会"
$wsAll.Range("D22").Value = "湖墅南路136-138号 浙话艺术剧院"
$wsAll.Range("E22").Value = "2024.07.13 19:30-07.13 21:30"
$wsAll.Range("F22").Value = 12
$wsAll.Range("G22").Value = 60
$wsAll.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=85011"
$wsAll.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202404/2Gd8eLva1714379746993.jpeg"

$wsAll.Range("B23").Value = "2024-07-13"
$wsAll.Range("C23").Value = "杭州·代号鸢only-广陵大学"
$wsAll.Range("D23").Value = "康候圣街99号 顺丰创新中心"
$wsAll.Range("E23").Value = "2024.07.13 09:00-07.13 18:00"
$wsAll.Range("F23").Value = 1165
$wsAll.Range("G23").Value = 68
$wsAll.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=83289"
$wsAll.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202403/I3yffJ7Q1711344958258.png"

$wsAll.Range("B24").Value = "2024-07-13"
$wsAll.Range("C24").Value = "杭州·海上钢琴师—一生必听的电影名曲《泰坦尼克号》《花样年华》《海上钢琴师》"
$wsAll.Range("D24").Value = "曙光路31号 浙江音乐厅"
$wsAll.Range("E24").Value = "2024.07.13 19:30-07.13 21:00"
$wsAll.Range("F24").Value = 3
$wsAll.Range("G24").Value = 100
$wsAll.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=85889"
$wsAll.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202405/52kxbBTh1716096935602.jpeg"

$wsAll.Range("B25").Value = "2024-07-13"
$wsAll.Range("C25").Value = "杭州·草莓动漫节"
$wsAll.Range("D25").Value = "中心路1号 白蓝地文创街区"
$wsAll.Range("E25").Value = "2024.07.13 09:00-07.14 17:00"
$wsAll.Range("F25").Value = 2903
$wsAll.Range("G25").Value = 70
$wsAll.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=84229"
$wsAll.Range("I25").Value = "//i2.hdslb.com/bfs/openplatform/202405/yjYrwO301715760081303.jpeg"

$wsAll.Range("B26").Value = "2024-07-20"
$wsAll.Range("C26").Value = "【会员购严选】杭州·首届次元格子动漫展-进入格子空间，探索次元世界！"
$wsAll.Range("D26").Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$wsAll.Range("E26").Value = "2024.07.20 09:00-07.22 17:00"
$wsAll.Range("F26").Value = 2260
$wsAll.Range("G26").Value = 75
$wsAll.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=85616"
$wsAll.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202405/5Dne5VqI1715753018080.jpeg"

$wsAll.Range("B27").Value = "2024-07-20"
$wsAll.Range("C27").Value = "杭州·TCD国潮动漫游戏嘉年华"
$wsAll.Range("D27").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$wsAll.Range("E27").Value = "2024.07.20 09:30-07.21 17:00"
$wsAll.Range("F27").Value = 4207
$wsAll.Range("G27").Value = 65
$wsAll.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=85699"
$wsAll.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202405/SwLIIdfR1715251191803.jpeg"

